$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.994.73'
$ws.Range('E2').Value = '  +2.08%  '
$ws.Range('D3').Value = '3.186.84'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''535.16'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '''145.68'
$ws.Range('E6').Value = '  +3.74%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''0.526'
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('D11').Value = '''0.430'
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').Value = '3.733.55'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = '''25.88'
$ws.Range('E14').Value = '  -1.34%  '
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '59.990.73'
$ws.Range('E16').Value = '  +1.99%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.197.55'
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '''6.27'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = '''13.23'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').Value = '''8.20'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '''368.58'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').Value = '''0.521'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '''69.46'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').Value = '0.0₃0869'
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').Value = '''22.34'
$ws.Range('E29').Value = '  +0.82%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '''6.05'
$ws.Range('E31').Value = '  -1.22%  '
$ws.Range('E32').Value = '  +1.56%  '
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').Value = '''6.56'
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('D35').Value = '''155.82'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').Value = '2.814.29'
$ws.Range('E37').Value = '  +7.01%  '
$ws.Range('D38').Value = '''26.12'
$ws.Range('E38').Value = '  +3.48%  '
$ws.Range('E39').Value = '  +2.73%  '
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').Value = '''39.69'
$ws.Range('E43').Value = '  +2.33%  '
$ws.Range('D44').Value = '''0.717'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').Value = '3.227.01'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').Value = '''0.984'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').Value = '''20.70'
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').Value = '''0.793'
$ws.Range('E50').Value = '  +4.46%  '
$ws.Range('D51').Value = '''0.999'
$ws.Range('E51').Value = '  -0.06%  '
